$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1010.2727  # H28
$ws.Cells.Item(28, 9).Value = 123.75  # I28
$ws.Cells.Item(28, 10).Value = 3374.3333  # J28
$ws.Cells.Item(28, 11).Value = 123.75  # K28
$ws.Cells.Item(28, 12).Value = 3374.3333  # L28
$ws.Cells.Item(28, 13).Value = 361.25  # M28
$ws.Cells.Item(28, 14).Value = -4344.3333  # N28
$ws.Cells.Item(43, 8).Value = 83334650  # H43
$ws.Cells.Item(43, 9).Value = 200000580  # I43
$ws.Cells.Item(43, 10).Value = 18520236  # J43
$ws.Cells.Item(43, 11).Value = 200000580  # K43
$ws.Cells.Item(43, 12).Value = 18520236  # L43
$ws.Cells.Item(43, 13).Value = -200000511  # M43
$ws.Cells.Item(43, 14).Value = -18520374  # N43
$ws.Cells.Item(113, 8).Value = 2747.5925  # H113
$ws.Cells.Item(113, 9).Value = 2086.818  # I113
$ws.Cells.Item(113, 10).Value = 3201.875  # J113
$ws.Cells.Item(113, 11).Value = 2086.818  # K113
$ws.Cells.Item(113, 12).Value = 3201.875  # L113
$ws.Cells.Item(113, 13).Value = 1167.182  # M113
$ws.Cells.Item(113, 14).Value = -9709.875  # N113
$ws.Cells.Item(116, 8).Value = 2062.5  # H116
$ws.Cells.Item(116, 9).Value = 2000  # I116
$ws.Cells.Item(116, 10).Value = 2500  # J116
$ws.Cells.Item(116, 11).Value = 2000  # K116
$ws.Cells.Item(116, 12).Value = 2500  # L116
$ws.Cells.Item(116, 13).Value = 1442  # M116
$ws.Cells.Item(116, 14).Value = -9384  # N116
$ws.Cells.Item(125, 8).Value = 897.087  # H125
$ws.Cells.Item(125, 9).Value = 475  # I125
$ws.Cells.Item(125, 10).Value = 1081.75  # J125
$ws.Cells.Item(125, 11).Value = 4275  # K125
$ws.Cells.Item(125, 12).Value = 9735.75  # L125
$ws.Cells.Item(125, 13).Value = -1815  # M125
$ws.Cells.Item(125, 14).Value = -14655.75  # N125
$ws.Cells.Item(132, 8).Value = 2917485.2  # H132
$ws.Cells.Item(132, 9).Value = 3970430.2  # I132
$ws.Cells.Item(132, 10).Value = 1637.5385  # J132
$ws.Cells.Item(132, 11).Value = 11911290.6  # K132
$ws.Cells.Item(132, 12).Value = 4912.6155  # L132
$ws.Cells.Item(132, 13).Value = -11908760.6  # M132
$ws.Cells.Item(132, 14).Value = -9972.6155  # N132
$ws.Cells.Item(135, 8).Value = 613.5574  # H135
$ws.Cells.Item(135, 9).Value = 429.77585  # I135
$ws.Cells.Item(135, 11).Value = 3867.98265  # K135
$ws.Cells.Item(135, 13).Value = -1332.98265  # M135
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 657.6799999999999  # H32
$ws.Cells.Item(32, 9).Value = 654.2222  # I32
$ws.Cells.Item(32, 11).Value = 654.2222  # K32
$ws.Cells.Item(32, 13).Value = -367.2222  # M32
$ws.Cells.Item(45, 8).Value = 1268.12  # H45
$ws.Cells.Item(45, 9).Value = 1098.091  # I45
$ws.Cells.Item(45, 10).Value = 1401.7142  # J45
$ws.Cells.Item(45, 11).Value = 1098.091  # K45
$ws.Cells.Item(45, 12).Value = 1401.7142  # L45
$ws.Cells.Item(45, 13).Value = -721.0909999999999  # M45
$ws.Cells.Item(45, 14).Value = -2155.7142  # N45
$ws.Cells.Item(61, 8).Value = 1631.3422  # H61
$ws.Cells.Item(61, 9).Value = 1018.3461  # I61
$ws.Cells.Item(61, 10).Value = 2959.5  # J61
$ws.Cells.Item(61, 11).Value = 1018.3461  # K61
$ws.Cells.Item(61, 12).Value = 2959.5  # L61
$ws.Cells.Item(61, 13).Value = -806.3461  # M61
$ws.Cells.Item(61, 14).Value = -3383.5  # N61
$ws.Cells.Item(74, 8).Value = 688.913  # H74
$ws.Cells.Item(74, 9).Value = 578.8372000000001  # I74
$ws.Cells.Item(74, 10).Value = 2266.6667  # J74
$ws.Cells.Item(74, 11).Value = 578.8372000000001  # K74
$ws.Cells.Item(74, 12).Value = 2266.6667  # L74
$ws.Cells.Item(74, 13).Value = 295.1627999999999  # M74
$ws.Cells.Item(74, 14).Value = -4014.6667  # N74
$ws.Cells.Item(76, 8).Value = 0  # H76
$ws.Cells.Item(76, 10).Value = 0  # J76
$ws.Cells.Item(76, 12).Value = 0  # L76
$ws.Cells.Item(76, 14).ClearContents()  # N76 (removed)
$ws.Cells.Item(77, 8).Value = 688.913  # H77
$ws.Cells.Item(77, 9).Value = 578.8372000000001  # I77
$ws.Cells.Item(77, 10).Value = 2266.6667  # J77
$ws.Cells.Item(77, 11).Value = 2894.186  # K77
$ws.Cells.Item(77, 12).Value = 11333.3335  # L77
$ws.Cells.Item(77, 13).Value = 1473.814  # M77
$ws.Cells.Item(77, 14).Value = -20069.3335  # N77
$ws.Cells.Item(79, 8).Value = 0  # H79
$ws.Cells.Item(79, 10).Value = 0  # J79
$ws.Cells.Item(79, 12).Value = 0  # L79
$ws.Cells.Item(79, 14).ClearContents()  # N79 (removed)
$ws.Cells.Item(97, 8).Value = 421.08572  # H97
$ws.Cells.Item(97, 9).Value = 399.26666  # I97
$ws.Cells.Item(97, 10).Value = 552  # J97
$ws.Cells.Item(97, 11).Value = 399.26666  # K97
$ws.Cells.Item(97, 12).Value = 552  # L97
$ws.Cells.Item(97, 13).Value = 96.73334  # M97
$ws.Cells.Item(97, 14).Value = -1544  # N97
$ws.Cells.Item(122, 8).Value = 1151.7142  # H122
$ws.Cells.Item(122, 9).Value = 1046.32  # I122
$ws.Cells.Item(122, 10).Value = 2030  # J122
$ws.Cells.Item(122, 11).Value = 3138.96  # K122
$ws.Cells.Item(122, 12).Value = 6090  # L122
$ws.Cells.Item(122, 13).Value = -688.96  # M122
$ws.Cells.Item(122, 14).Value = -10990  # N122
$ws.Cells.Item(132, 8).Value = 5084.237  # H132
$ws.Cells.Item(132, 9).Value = 5803.8076  # I132
$ws.Cells.Item(132, 10).Value = 3525.1667  # J132
$ws.Cells.Item(132, 11).Value = 17411.4228  # K132
$ws.Cells.Item(132, 12).Value = 10575.5001  # L132
$ws.Cells.Item(132, 13).Value = -14881.4228  # M132
$ws.Cells.Item(132, 14).Value = -15635.5001  # N132
$ws.Cells.Item(136, 8).Value = 1631.3422  # H136
$ws.Cells.Item(136, 9).Value = 1018.3461  # I136
$ws.Cells.Item(136, 10).Value = 2959.5  # J136
$ws.Cells.Item(136, 11).Value = 3055.0383  # K136
$ws.Cells.Item(136, 12).Value = 8878.5  # L136
$ws.Cells.Item(136, 13).Value = -505.0383000000002  # M136
$ws.Cells.Item(136, 14).Value = -13978.5  # N136
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 24165.268  # H134
$ws.Cells.Item(134, 9).Value = 33937.13  # I134
$ws.Cells.Item(134, 10).Value = 2527.5715  # J134
$ws.Cells.Item(134, 11).Value = 101811.39  # K134
$ws.Cells.Item(134, 12).Value = 7582.7145  # L134
$ws.Cells.Item(134, 13).Value = -99276.38999999998  # M134
$ws.Cells.Item(134, 14).Value = -12652.7145  # N134
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3032243  # H31
$ws.Cells.Item(31, 9).Value = 1442.4147  # I31
$ws.Cells.Item(31, 10).Value = 11908159  # J31
$ws.Cells.Item(31, 11).Value = 1442.4147  # K31
$ws.Cells.Item(31, 12).Value = 11908159  # L31
$ws.Cells.Item(31, 13).Value = -1147.4147  # M31
$ws.Cells.Item(31, 14).Value = -11908749  # N31
$ws.Cells.Item(34, 8).Value = 3032243  # H34
$ws.Cells.Item(34, 9).Value = 1442.4147  # I34
$ws.Cells.Item(34, 10).Value = 11908159  # J34
$ws.Cells.Item(34, 11).Value = 1442.4147  # K34
$ws.Cells.Item(34, 12).Value = 11908159  # L34
$ws.Cells.Item(34, 13).Value = -1240.4147  # M34
$ws.Cells.Item(34, 14).Value = -11908563  # N34
$ws.Cells.Item(58, 8).Value = 5556204  # H58
$ws.Cells.Item(58, 9).Value = 674.31915  # I58
$ws.Cells.Item(58, 10).Value = 25641580  # J58
$ws.Cells.Item(58, 11).Value = 674.31915  # K58
$ws.Cells.Item(58, 12).Value = 25641580  # L58
$ws.Cells.Item(58, 13).Value = -471.31915  # M58
$ws.Cells.Item(58, 14).Value = -25641986  # N58
$ws.Cells.Item(132, 8).Value = 1833.6531  # H132
$ws.Cells.Item(132, 9).Value = 2023.2188  # I132
$ws.Cells.Item(132, 10).Value = 1476.8235  # J132
$ws.Cells.Item(132, 11).Value = 6069.6564  # K132
$ws.Cells.Item(132, 12).Value = 4430.470499999999  # L132
$ws.Cells.Item(132, 13).Value = -3539.6564  # M132
$ws.Cells.Item(132, 14).Value = -9490.470499999999  # N132
$ws.Cells.Item(134, 8).Value = 933.14923  # H134
$ws.Cells.Item(134, 9).Value = 908.3090999999999  # I134
$ws.Cells.Item(134, 10).Value = 1047  # J134
$ws.Cells.Item(134, 11).Value = 2724.9273  # K134
$ws.Cells.Item(134, 12).Value = 3141  # L134
$ws.Cells.Item(134, 13).Value = -189.9272999999998  # M134
$ws.Cells.Item(134, 14).Value = -8211  # N134
$ws.Cells.Item(136, 8).Value = 5556204  # H136
$ws.Cells.Item(136, 9).Value = 674.31915  # I136
$ws.Cells.Item(136, 10).Value = 25641580  # J136
$ws.Cells.Item(136, 11).Value = 2022.95745  # K136
$ws.Cells.Item(136, 12).Value = 76924740  # L136
$ws.Cells.Item(136, 13).Value = 527.0425499999999  # M136
$ws.Cells.Item(136, 14).Value = -76929840  # N136
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 567.1429000000001  # H46
$ws.Cells.Item(46, 10).Value = 842.5  # J46
$ws.Cells.Item(46, 12).Value = 2527.5  # L46
$ws.Cells.Item(46, 14).Value = -2709.5  # N46
$ws.Cells.Item(116, 8).Value = 9631.5  # H116
$ws.Cells.Item(116, 9).Value = 11953.111  # I116
$ws.Cells.Item(116, 11).Value = 35859.333  # K116
$ws.Cells.Item(116, 13).Value = -32417.333  # M116
$ws.Cells.Item(131, 8).Value = 1864725.4  # H131
$ws.Cells.Item(131, 9).Value = 5684.7144  # I131
$ws.Cells.Item(131, 10).Value = 2587685.8  # J131
$ws.Cells.Item(131, 11).Value = 17054.1432  # K131
$ws.Cells.Item(131, 12).Value = 7763057.399999999  # L131
$ws.Cells.Item(131, 13).Value = -12014.1432  # M131
$ws.Cells.Item(131, 14).Value = -7773137.399999999  # N131
# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(95, 8).Value = 7983.3335  # H95
$ws.Cells.Item(95, 10).Value = 7983.3335  # J95
$ws.Cells.Item(95, 12).Value = 7983.3335  # L95
$ws.Cells.Item(95, 14).Value = -13475.3335  # N95
$ws.Cells.Item(102, 8).Value = 1511.5555  # H102
$ws.Cells.Item(102, 9).Value = 1550  # I102
$ws.Cells.Item(102, 11).Value = 1550  # K102
$ws.Cells.Item(102, 13).Value = 72  # M102
$ws.Cells.Item(126, 8).Value = 1888  # H126
$ws.Cells.Item(126, 9).Value = 1868.5  # I126
$ws.Cells.Item(126, 10).Value = 1940  # J126
$ws.Cells.Item(126, 11).Value = 5605.5  # K126
$ws.Cells.Item(126, 12).Value = 5820  # L126
$ws.Cells.Item(126, 13).Value = -3135.5  # M126
$ws.Cells.Item(126, 14).Value = -10760  # N126
$ws.Cells.Item(132, 8).Value = 31384.416  # H132
$ws.Cells.Item(132, 9).Value = 42174.8  # I132
$ws.Cells.Item(132, 10).Value = 6860.8184  # J132
$ws.Cells.Item(132, 11).Value = 126524.4  # K132
$ws.Cells.Item(132, 12).Value = 20582.4552  # L132
$ws.Cells.Item(132, 13).Value = -123994.4  # M132
$ws.Cells.Item(132, 14).Value = -25642.4552  # N132
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3250  # H22
$ws.Cells.Item(22, 9).Value = 4000  # I22
$ws.Cells.Item(22, 10).Value = 1000  # J22
$ws.Cells.Item(22, 11).Value = 4000  # K22
$ws.Cells.Item(22, 12).Value = 1000  # L22
$ws.Cells.Item(22, 13).Value = -3705  # M22
$ws.Cells.Item(22, 14).Value = -1590  # N22
$ws.Cells.Item(27, 8).Value = 3250  # H27
$ws.Cells.Item(27, 9).Value = 4000  # I27
$ws.Cells.Item(27, 10).Value = 1000  # J27
$ws.Cells.Item(27, 11).Value = 4000  # K27
$ws.Cells.Item(27, 12).Value = 1000  # L27
$ws.Cells.Item(27, 13).Value = -3893  # M27
$ws.Cells.Item(27, 14).Value = -1214  # N27
$ws.Cells.Item(61, 8).Value = 1835.0714  # H61
$ws.Cells.Item(61, 9).Value = 1650.1  # I61
$ws.Cells.Item(61, 10).Value = 2297.5  # J61
$ws.Cells.Item(61, 11).Value = 1650.1  # K61
$ws.Cells.Item(61, 12).Value = 2297.5  # L61
$ws.Cells.Item(61, 13).Value = -1448.1  # M61
$ws.Cells.Item(61, 14).Value = -2701.5  # N61
$ws.Cells.Item(106, 8).Value = 20207.777  # H106
$ws.Cells.Item(106, 10).Value = 20207.777  # J106
$ws.Cells.Item(106, 12).Value = 20207.777  # L106
$ws.Cells.Item(106, 14).Value = -22731.777  # N106
$ws.Cells.Item(113, 8).Value = 1835.0714  # H113
$ws.Cells.Item(113, 9).Value = 1650.1  # I113
$ws.Cells.Item(113, 10).Value = 2297.5  # J113
$ws.Cells.Item(113, 11).Value = 1650.1  # K113
$ws.Cells.Item(113, 12).Value = 2297.5  # L113
$ws.Cells.Item(113, 13).Value = 519.9000000000001  # M113
$ws.Cells.Item(113, 14).Value = -6637.5  # N113
$ws.Cells.Item(122, 8).Value = 2884.111  # H122
$ws.Cells.Item(122, 9).Value = 3383.2727  # I122
$ws.Cells.Item(122, 10).Value = 2540.9375  # J122
$ws.Cells.Item(122, 11).Value = 10149.8181  # K122
$ws.Cells.Item(122, 12).Value = 7622.8125  # L122
$ws.Cells.Item(122, 13).Value = -7699.8181  # M122
$ws.Cells.Item(122, 14).Value = -12522.8125  # N122
$ws.Cells.Item(132, 8).Value = 4218.05  # H132
$ws.Cells.Item(132, 9).Value = 5053.4443  # I132
$ws.Cells.Item(132, 10).Value = 1711.8667  # J132
$ws.Cells.Item(132, 11).Value = 15160.3329  # K132
$ws.Cells.Item(132, 12).Value = 5135.6001  # L132
$ws.Cells.Item(132, 13).Value = -12630.3329  # M132
$ws.Cells.Item(132, 14).Value = -10195.6001  # N132
$ws.Cells.Item(136, 8).Value = 2093.842  # H136
$ws.Cells.Item(136, 9).Value = 2086.1562  # I136
$ws.Cells.Item(136, 10).Value = 2134.8333  # J136
$ws.Cells.Item(136, 11).Value = 6258.4686  # K136
$ws.Cells.Item(136, 12).Value = 6404.499899999999  # L136
$ws.Cells.Item(136, 13).Value = -3708.4686  # M136
$ws.Cells.Item(136, 14).Value = -11504.4999  # N136
# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 11371.091  # H126
$ws.Cells.Item(126, 9).Value = 13788.777  # I126
$ws.Cells.Item(126, 10).Value = 491.5  # J126
$ws.Cells.Item(126, 11).Value = 41366.331  # K126
$ws.Cells.Item(126, 12).Value = 1474.5  # L126
$ws.Cells.Item(126, 13).Value = -38896.331  # M126
$ws.Cells.Item(126, 14).Value = -6414.5  # N126
$ws.Cells.Item(132, 8).Value = 1134.6383  # H132
$ws.Cells.Item(132, 9).Value = 1083.881  # I132
$ws.Cells.Item(132, 10).Value = 1561  # J132
$ws.Cells.Item(132, 11).Value = 3251.643  # K132
$ws.Cells.Item(132, 12).Value = 4683  # L132
$ws.Cells.Item(132, 13).Value = -721.643  # M132
$ws.Cells.Item(132, 14).Value = -9743  # N132
$ws.Cells.Item(136, 8).Value = 2200.5  # H136
$ws.Cells.Item(136, 9).Value = 2440.625  # I136
$ws.Cells.Item(136, 10).Value = 1240  # J136
$ws.Cells.Item(136, 11).Value = 7321.875  # K136
$ws.Cells.Item(136, 12).Value = 3720  # L136
$ws.Cells.Item(136, 13).Value = -4771.875  # M136
$ws.Cells.Item(136, 14).Value = -8820  # N136
